$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 286.7143
$ws.Range("I9").Value = 126.25
$ws.Range("K9").Value = 126.25
$ws.Range("M9").Value = 42.75
$ws.Range("H12").Value = 1329.091
$ws.Range("I12").Value = 1746.25
$ws.Range("J12").Value = 216.66667
$ws.Range("K12").Value = 1746.25
$ws.Range("L12").Value = 216.66667
$ws.Range("M12").Value = -1576.25
$ws.Range("N12").Value = -556.6666700000001
$ws.Range("H87").Value = 39854
$ws.Range("J87").Value = 39854
$ws.Range("L87").Value = 39854
$ws.Range("N87").Value = -42350
$ws.Range("H90").Value = 39854
$ws.Range("J90").Value = 39854
$ws.Range("L90").Value = 119562
$ws.Range("N90").Value = -132042
$ws.Range("H138").Value = 3942.3794
$ws.Range("I138").Value = 7166.3335
$ws.Range("J138").Value = 3766.5273
$ws.Range("K138").Value = 21499.0005
$ws.Range("L138").Value = 11299.5819
$ws.Range("M138").Value = -16359.0005
$ws.Range("N138").Value = -21579.5819

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2225.7896
$ws.Range("I61").Value = 1793.2188
$ws.Range("J61").Value = 4532.8335
$ws.Range("K61").Value = 1793.2188
$ws.Range("L61").Value = 4532.8335
$ws.Range("M61").Value = -1581.2188
$ws.Range("N61").Value = -4956.8335
$ws.Range("H74").Value = 21740460
$ws.Range("I74").Value = 28572014
$ws.Range("J74").Value = 3702.5454
$ws.Range("K74").Value = 28572014
$ws.Range("L74").Value = 3702.5454
$ws.Range("M74").Value = -28571140
$ws.Range("N74").Value = -5450.5454
$ws.Range("H77").Value = 21740460
$ws.Range("I77").Value = 28572014
$ws.Range("J77").Value = 3702.5454
$ws.Range("K77").Value = 142860070
$ws.Range("L77").Value = 18512.727
$ws.Range("M77").Value = -142855702
$ws.Range("N77").Value = -27248.727
$ws.Range("H80").Value = 50133.75
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 50133.75
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 50133.75
$ws.Range("M80").Value = ""
$ws.Range("N80").Value = -52129.75
$ws.Range("H83").Value = 50133.75
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 50133.75
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 150401.25
$ws.Range("M83").Value = ""
$ws.Range("N83").Value = -160385.25
$ws.Range("H97").Value = 763.6
$ws.Range("I97").Value = 773.1111
$ws.Range("J97").Value = 678
$ws.Range("K97").Value = 773.1111
$ws.Range("L97").Value = 678
$ws.Range("M97").Value = -277.1111
$ws.Range("N97").Value = -1670
$ws.Range("H132").Value = 14810
$ws.Range("I132").Value = 2291.4
$ws.Range("J132").Value = 48951.637
$ws.Range("K132").Value = 6874.200000000001
$ws.Range("L132").Value = 146854.911
$ws.Range("M132").Value = -4344.200000000001
$ws.Range("N132").Value = -151914.911
$ws.Range("H136").Value = 2225.7896
$ws.Range("I136").Value = 1793.2188
$ws.Range("J136").Value = 4532.8335
$ws.Range("K136").Value = 5379.6564
$ws.Range("L136").Value = 13598.5005
$ws.Range("M136").Value = -2829.6564
$ws.Range("N136").Value = -18698.5005

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 35965.332
$ws.Range("I82").Value = 8257
$ws.Range("J82").Value = 41507
$ws.Range("K82").Value = 8257
$ws.Range("L82").Value = 41507
$ws.Range("M82").Value = -7874
$ws.Range("N82").Value = -42273
$ws.Range("H85").Value = 35965.332
$ws.Range("I85").Value = 8257
$ws.Range("J85").Value = 41507
$ws.Range("K85").Value = 8257
$ws.Range("L85").Value = 41507
$ws.Range("M85").Value = -6931
$ws.Range("N85").Value = -44159
$ws.Range("H99").Value = 1410.1333
$ws.Range("I99").Value = 1316.6666
$ws.Range("K99").Value = 1316.6666
$ws.Range("M99").Value = 181.3334
$ws.Range("H134").Value = 3302.9333
$ws.Range("I134").Value = 3469.476
$ws.Range("J134").Value = 971.3333
$ws.Range("K134").Value = 10408.428
$ws.Range("L134").Value = 2913.9999
$ws.Range("M134").Value = -7873.428
$ws.Range("N134").Value = -7983.9999
$ws.Range("H140").Value = 30780
$ws.Range("J140").Value = 30780
$ws.Range("L140").Value = 30780
$ws.Range("N140").Value = -41140

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3151.2568
$ws.Range("I31").Value = 1510.3
$ws.Range("J31").Value = 5081.794
$ws.Range("K31").Value = 1510.3
$ws.Range("L31").Value = 5081.794
$ws.Range("M31").Value = -1215.3
$ws.Range("N31").Value = -5671.794
$ws.Range("H34").Value = 3151.2568
$ws.Range("I34").Value = 1510.3
$ws.Range("J34").Value = 5081.794
$ws.Range("K34").Value = 1510.3
$ws.Range("L34").Value = 5081.794
$ws.Range("M34").Value = -1308.3
$ws.Range("N34").Value = -5485.794
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").Value = ""
$ws.Range("H58").Value = 20837.846
$ws.Range("I58").Value = 1601.0714
$ws.Range("J58").Value = 43280.75
$ws.Range("K58").Value = 1601.0714
$ws.Range("L58").Value = 43280.75
$ws.Range("M58").Value = -1398.0714
$ws.Range("N58").Value = -43686.75
$ws.Range("H132").Value = 2570.2415
$ws.Range("I132").Value = 1814.3334
$ws.Range("J132").Value = 4554.5
$ws.Range("K132").Value = 5443.0002
$ws.Range("L132").Value = 13663.5
$ws.Range("M132").Value = -2913.0002
$ws.Range("N132").Value = -18723.5
$ws.Range("H134").Value = 1252.6111
$ws.Range("I134").Value = 1031.8572
$ws.Range("J134").Value = 2025.25
$ws.Range("K134").Value = 3095.5716
$ws.Range("L134").Value = 6075.75
$ws.Range("M134").Value = -560.5715999999998
$ws.Range("N134").Value = -11145.75
$ws.Range("H136").Value = 20837.846
$ws.Range("I136").Value = 1601.0714
$ws.Range("J136").Value = 43280.75
$ws.Range("K136").Value = 4803.2142
$ws.Range("L136").Value = 129842.25
$ws.Range("M136").Value = -2253.2142
$ws.Range("N136").Value = -134942.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 25595.5
$ws.Range("J68").Value = 50501.5
$ws.Range("L68").Value = 151504.5
$ws.Range("N68").Value = -153126.5
$ws.Range("H69").Value = 2433.3333
$ws.Range("J69").Value = 2450
$ws.Range("L69").Value = 7350
$ws.Range("N69").Value = -8972
$ws.Range("H71").Value = 25595.5
$ws.Range("J71").Value = 50501.5
$ws.Range("L71").Value = 454513.5
$ws.Range("N71").Value = -462625.5
$ws.Range("H72").Value = 2433.3333
$ws.Range("J72").Value = 2450
$ws.Range("L72").Value = 22050
$ws.Range("N72").Value = -30162
$ws.Range("H92").Value = 25000770
$ws.Range("J92").Value = 1499.5
$ws.Range("L92").Value = 4498.5
$ws.Range("N92").Value = -6994.5
$ws.Range("H131").Value = 734.63
$ws.Range("J131").Value = 764.43475
$ws.Range("L131").Value = 2293.30425
$ws.Range("N131").Value = -12373.30425
$ws.Range("H132").Value = 1443.3334
$ws.Range("I132").Value = 1415
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 12735
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -10205
$ws.Range("N132").Value = -18560

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1718.5
$ws.Range("I102").Value = 1778.3334
$ws.Range("K102").Value = 1778.3334
$ws.Range("M102").Value = -156.3334
$ws.Range("H122").Value = 4126.8276
$ws.Range("I122").Value = 3983.1052
$ws.Range("J122").Value = 4399.9
$ws.Range("K122").Value = 11949.3156
$ws.Range("L122").Value = 13199.7
$ws.Range("M122").Value = -9499.3156
$ws.Range("N122").Value = -18099.7
$ws.Range("H132").Value = 22357.178
$ws.Range("I132").Value = 5043.8335
$ws.Range("J132").Value = 53521.2
$ws.Range("K132").Value = 15131.5005
$ws.Range("L132").Value = 160563.6
$ws.Range("M132").Value = -12601.5005
$ws.Range("N132").Value = -165623.6

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I22").Value = 4380.4
$ws.Range("J22").Value = 1386
$ws.Range("K22").Value = 4380.4
$ws.Range("L22").Value = 1386
$ws.Range("M22").Value = -4085.4
$ws.Range("N22").Value = -1976
$ws.Range("I27").Value = 4380.4
$ws.Range("J27").Value = 1386
$ws.Range("K27").Value = 4380.4
$ws.Range("L27").Value = 1386
$ws.Range("M27").Value = -4273.4
$ws.Range("N27").Value = -1600
$ws.Range("H46").Value = 1277.8889
$ws.Range("I46").Value = 966.8333
$ws.Range("J46").Value = 1900
$ws.Range("K46").Value = 966.8333
$ws.Range("L46").Value = 1900
$ws.Range("M46").Value = -778.8333
$ws.Range("N46").Value = -2276
$ws.Range("H82").Value = 1490.6666
$ws.Range("I82").Value = 1490
$ws.Range("K82").Value = 1490
$ws.Range("M82").Value = -1129
$ws.Range("H85").Value = 1490.6666
$ws.Range("I85").Value = 1490
$ws.Range("K85").Value = 1490
$ws.Range("M85").Value = -242
$ws.Range("H132").Value = 243647.9
$ws.Range("I132").Value = 356688.12
$ws.Range("J132").Value = 3437.4375
$ws.Range("K132").Value = 1070064.36
$ws.Range("L132").Value = 10312.3125
$ws.Range("M132").Value = -1067534.36
$ws.Range("N132").Value = -15372.3125
$ws.Range("H136").Value = 2122.0588
$ws.Range("I136").Value = 1957.8125
$ws.Range("K136").Value = 5873.4375
$ws.Range("M136").Value = -3323.4375

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 12165.833
$ws.Range("I41").Value = 5500
$ws.Range("J41").Value = 13499
$ws.Range("K41").Value = 5500
$ws.Range("L41").Value = 13499
$ws.Range("M41").Value = -5110
$ws.Range("N41").Value = -14279
$ws.Range("H75").Value = 24000
$ws.Range("J75").Value = 24000
$ws.Range("L75").Value = 24000
$ws.Range("N75").Value = -25872
$ws.Range("H78").Value = 24000
$ws.Range("J78").Value = 24000
$ws.Range("L78").Value = 72000
$ws.Range("N78").Value = -81360
$ws.Range("H107").Value = 200000460
$ws.Range("I107").Value = 250000400
$ws.Range("K107").Value = 750001200
$ws.Range("M107").Value = -749999280
$ws.Range("H136").Value = 28676820
$ws.Range("I136").Value = 35596540
$ws.Range("K136").Value = 106789620
$ws.Range("M136").Value = -106787070
